$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (not number/date) interpretation for the Price/Volume columns so that
# values like "0.999" or "316.78" are stored as text, matching the source data,
# then restore the original (default) cell style so no visible formatting changes.
$priceVolRange = $ws.Range("D2:E51")
$origStyle = $priceVolRange.Style
$priceVolRange.NumberFormat = "@"

$ws.Range('D2').Value = '47.797.30'
$ws.Range('D3').Value = '2.482.18'
$ws.Range('D5').Value = '316.78'
$ws.Range('D6').Value = '103.97'
$ws.Range('D7').Value = '0.517'
$ws.Range('D8').Value = '0.999'
$ws.Range('D9').Value = '0.533'
$ws.Range('D10').Value = '38.60'
$ws.Range('D11').Value = '20.44'
$ws.Range('D12').Value = '0.0797'
$ws.Range('D15').Value = '2.867.62'
$ws.Range('D16').Value = '2.511.33'
$ws.Range('D17').Value = '0.822'
$ws.Range('D18').Value = '47.716.84'
$ws.Range('D19').Value = '2.91'
$ws.Range('D21').Value = '6.51'
$ws.Range('D22').Value = '0.0₃0924'
$ws.Range('D23').Value = '277.75'
$ws.Range('D24').Value = '70.67'
$ws.Range('D27').Value = '25.63'
$ws.Range('D28').Value = '2.16'
$ws.Range('D29').Value = '9.56'
$ws.Range('D30').Value = '0.136'
$ws.Range('D31').Value = '34.36'
$ws.Range('D34').Value = '18.88'
$ws.Range('D35').Value = '5.24'
$ws.Range('D36').Value = '0.0766'
$ws.Range('D38').Value = '4.49'
$ws.Range('D40').Value = '122.16'
$ws.Range('D42').Value = '2.20'
$ws.Range('D43').Value = '21.54'
$ws.Range('D44').Value = '0.0298'
$ws.Range('D45').Value = '1.989.59'
$ws.Range('D46').Value = '3.12'
$ws.Range('D47').Value = '1.88'
$ws.Range('D49').Value = '8.87'
$ws.Range('D50').Value = '5.06'
$ws.Range('D51').Value = '78.49'

$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('E3').Value = '  -1.74%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('E5').Value = '  -1.74%  '
$ws.Range('E6').Value = '  -4.95%  '
$ws.Range('E7').Value = '  -2.92%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -3.62%  '
$ws.Range('E10').Value = '  -4.69%  '
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  -3.40%  '
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('E14').Value = '  -3.89%  '
$ws.Range('E15').Value = '  -1.78%  '
$ws.Range('E16').Value = '  -0.75%  '
$ws.Range('E17').Value = '  -3.86%  '
$ws.Range('E18').Value = '  -0.95%  '
$ws.Range('E19').Value = '  +7.71%  '
$ws.Range('E20').Value = '  -6.54%  '
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('E22').Value = '  -2.60%  '
$ws.Range('E23').Value = '  +4.94%  '
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('E25').Value = '  -3.48%  '
$ws.Range('E26').Value = '  -0.19%  '
$ws.Range('E27').Value = '  -1.83%  '
$ws.Range('E28').Value = '  -7.86%  '
$ws.Range('E29').Value = '  -5.59%  '
$ws.Range('E30').Value = '  -5.35%  '
$ws.Range('E31').Value = '  -3.85%  '
$ws.Range('E32').Value = '  -1.12%  '
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('E34').Value = '  -4.69%  '
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('E36').Value = '  -2.87%  '
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('E38').Value = '  -4.80%  '
$ws.Range('E39').Value = '  -5.46%  '
$ws.Range('E40').Value = '  +1.05%  '
$ws.Range('E41').Value = '  -1.80%  '
$ws.Range('E43').Value = '  -2.64%  '
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('E47').Value = '  -1.31%  '
$ws.Range('E48').Value = '  -4.19%  '
$ws.Range('E49').Value = '  -3.07%  '
$ws.Range('E50').Value = '  -3.12%  '
$ws.Range('E51').Value = '  -1.03%  '

# Restore original style/number-format so formatting matches the source workbook.
$priceVolRange.Style = $origStyle
